# Trade #66 closed at 2026-02-17 08:49:04 - unknown UNKNOWN +0.000%
#
# Updates the live trading results workbook:
#  - Summary sheet: refresh aggregate capital / P&L / trade counters
#  - Strategy Status sheet: refresh the MarketMaking strategy row
#  - All Trades / MarketMaking sheets: append the newly closed trade (#66)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1199.68   # Current Capital
$wsSummary.Range("B4").Value = -0.32     # Total P&L $
$wsSummary.Range("B5").Value = -0.1      # Total P&L %
$wsSummary.Range("B6").Value = 66        # Total Trades
$wsSummary.Range("B7").Value = 27        # Winning Trades
$wsSummary.Range("B9").Value = 40.91     # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet (MarketMaking row)
# ---------------------------------------------------------------------
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C4").Value = 99.68000000000001   # Capital
$wsStatus.Range("D4").Value = 66                  # Trades
$wsStatus.Range("E4").Value = -0.32               # P&L $
$wsStatus.Range("F4").Value = -0.32               # P&L %
$wsStatus.Range("G4").Value = 40.91               # Win Rate %

# ---------------------------------------------------------------------
# Append trade #66 to both the "All Trades" and "MarketMaking" logs
# ---------------------------------------------------------------------
$newRow = 67

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Date/Time columns look like numbers to Excel's auto-detection, so
    # force them to Text first to keep them as literal strings (matching
    # every other row in the log) instead of being parsed into date
    # serials.
    $ws.Range("B" + $newRow + ":C" + $newRow).NumberFormat = "@"

    $ws.Range("A" + $newRow).Value = 66
    $ws.Range("B" + $newRow).Value = "2026-02-17"
    $ws.Range("C" + $newRow).Value = "08:48:58"
    $ws.Range("D" + $newRow).Value = "MarketMaking"
    $ws.Range("E" + $newRow).Value = "UP"
    $ws.Range("F" + $newRow).Value = 0.92
    $ws.Range("G" + $newRow).Value = 0.98
    $ws.Range("H" + $newRow).Value = "CLOSED"
    $ws.Range("I" + $newRow).Value = 6.5217
    $ws.Range("J" + $newRow).Value = 0.06
    $ws.Range("K" + $newRow).Value = 99.68000000000001
    $ws.Range("L" + $newRow).Value = 0
    $ws.Range("M" + $newRow).Value = 0
    $ws.Range("N" + $newRow).Value = 0.6
    $ws.Range("O" + $newRow).Value = "Normal spread capture: 19600 bps"
    $ws.Range("P" + $newRow).Value = "early_exit"
    $ws.Range("Q" + $newRow).Value = 0.14

    # Drop the temporary text format override again so the new cells end
    # up with the same (default) styling as the rest of the sheet.
    $ws.Range("B" + $newRow + ":C" + $newRow).ClearFormats()
}
